$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.237.22"

# Row 3
$ws.Range("D3").Value = "1.860.61"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'0.7059"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("D6").Value = "'242.32"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'0.07819"
$ws.Range("E8").Value = "  -2.63%  "

# Row 9
$ws.Range("D9").Value = "'0.3113"
$ws.Range("E9").Value = "  -0.53%  "

# Row 10
$ws.Range("D10").Value = "'24.30"
$ws.Range("E10").Value = "  -3.78%  "

# Row 11
$ws.Range("D11").Value = "'0.08004"
$ws.Range("E11").Value = "  -4.28%  "

# Row 12
$ws.Range("D12").Value = "1.884.33"
$ws.Range("E12").Value = "  -0.97%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.175"
$ws.Range("E13").Value = "  -1.36%  "

# Row 14
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'93.36"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15
$ws.Range("D15").Value = "'0.6953"
$ws.Range("E15").Value = "  -3.55%  "

# Row 16
$ws.Range("D16").Value = "'6.350"
$ws.Range("E16").Value = "  +0.84%  "

# Row 17
$ws.Range("D17").Value = "29.268.37"
$ws.Range("E17").Value = "  -0.44%  "

# Row 18
$ws.Range("D18").Value = "'0.000008296"
$ws.Range("E18").Value = "  -2.07%  "

# Row 19
$ws.Range("D19").Value = "'252.89"
$ws.Range("E19").Value = "  +4.81%  "

# Row 20
$ws.Range("D20").Value = "2.159.36"
$ws.Range("E20").Value = "  +1.33%  "

# Row 21
$ws.Range("D21").Value = "'13.09"
$ws.Range("E21").Value = "  -1.23%  "

# Row 22
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").Value = "'7.527"
$ws.Range("E23").Value = "  -4.25%  "

# Row 24
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("E25").Value = "  -1.99%  "

# Row 26
$ws.Range("E26").Value = "  -0.95%  "

# Row 27
$ws.Range("D27").Value = "'159.46"
$ws.Range("E27").Value = "  -2.57%  "

# Row 28
$ws.Range("D28").Value = "'18.72"
$ws.Range("E28").Value = "  +0.76%  "

# Row 29
$ws.Range("D29").Value = "'1.497"
$ws.Range("E29").Value = "  -0.62%  "

# Row 30
$ws.Range("D30").Value = "'4.267"
$ws.Range("E30").Value = "  -1.72%  "

# Row 31
$ws.Range("D31").Value = "'4.270"
$ws.Range("E31").Value = "  -3.43%  "

# Row 32
$ws.Range("D32").Value = "'1.208"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("D33").Value = "'0.05268"

# Row 34
$ws.Range("D34").Value = "'1.887"
$ws.Range("E34").Value = "  -3.43%  "

# Row 35
$ws.Range("D35").Value = "'0.7429"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36
$ws.Range("D36").Value = "'1.155"
$ws.Range("E36").Value = "  -2.31%  "

# Row 37
$ws.Range("E37").Value = "  +0.24%  "

# Row 38
$ws.Range("E38").Value = "  -1.44%  "

# Row 39
$ws.Range("D39").Value = "1.249.76"
$ws.Range("E39").Value = "  -2.88%  "

# Row 40
$ws.Range("D40").Value = "'2.737"
$ws.Range("E40").Value = "  -0.42%  "

# Row 41
$ws.Range("D41").Value = "'6.283"
$ws.Range("E41").Value = "  -4.52%  "

# Row 42
$ws.Range("D42").Value = "'0.9021"
$ws.Range("E42").Value = "  +1.00%  "

# Row 43
$ws.Range("D43").Value = "'111.01"

# Row 44
$ws.Range("D44").Value = "'71.65"
$ws.Range("E44").Value = "  -2.89%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("D46").Value = "2.062.91"
$ws.Range("E46").Value = "  +1.64%  "

# Row 47
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("D48").Value = "'0.5202"
$ws.Range("E48").Value = "  -0.32%  "

# Row 49
$ws.Range("D49").Value = "'1.783"
$ws.Range("E49").Value = "  -1.36%  "

# Row 50
$ws.Range("D50").Value = "'9.390"
$ws.Range("E50").Value = "  -1.19%  "

# Row 51
$ws.Range("D51").Value = "'1.009"
$ws.Range("E51").Value = "  +0.88%  "
